$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "H2"  = 3.25
    "I2"  = 4.33
    "K2"  = 1.95
    "L2"  = 5
    "O2"  = 1.5
    "P2"  = 2.5
    "Q2"  = 2.5
    "R2"  = 1.5
    "S2"  = 1.57
    "T2"  = 2.25
    "U2"  = 2.2
    "V2"  = 1.62
    "W2"  = 5.5
    "AB2" = 41
    "AC2" = 6.5
    "AE2" = 21
    "AF2" = 81
    "AL2" = 51
    "AT2" = 2.25
    "AU2" = 9.5
    "AV2" = 81
    "AZ2" = 101
    "BA2" = 151
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
